$d = $word.ActiveDocument

function Insert-RunsXml($Range, $InnerXml) {
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($xml)
}

# --- Paragraph 2: "{{name}}" -> "{{firstName}} {{lastName}}" (with proofErr spell-check runs) ---
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1
$p2Range = $d.Range($p2Start, $p2End)

$inner2 = '<w:r><w:t>{{</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>firstN</w:t></w:r>' +
          '<w:r><w:t>ame</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>}}</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:r><w:t>{{</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>last</w:t></w:r>' +
          '<w:r><w:t>Name</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>}}</w:t></w:r>'

Insert-RunsXml $p2Range $inner2

# --- Paragraph 3: "{{" + "#" + "end" + "each}}" -> "{{" + "#end" + "each}}" ---
$p3 = $d.Paragraphs.Item(3)
$p3Start = $p3.Range.Start
$p3End = $p3.Range.End - 1
$p3Range = $d.Range($p3Start, $p3End)

$inner3 = '<w:r><w:t>{{</w:t></w:r>' +
          '<w:r w:rsidR="00DD4E36"><w:t>#end</w:t></w:r>' +
          '<w:r><w:t>each}}</w:t></w:r>'

Insert-RunsXml $p3Range $inner3

Write-Host "Done"
